# SC List view on vendor panel, Downloadable in our CRM
# -> Remove the "State" / {sc:state} column from the SC Charges List template.
# Deleting the entire column shifts every later column (C..K) one slot to
# the left (C->B, D->C, ... K->J), which matches the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B:B").EntireColumn.Delete()

# Leave the same cell/column selected as in the authored workbook.
$ws.Range("B1:B1048576").Select()
